$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 = "I0" and J1 = "IF", matching the formatting of the
# existing header cell (e.g. H1) by copying its format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in data for rows 2-34: column I is always 1, column J mirrors
# the value already present in column H (the "IP" column) for that row.
for ($r = 2; $r -le 34; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
